# feat: added img url scanning feature in CSVtoPLANTS
#
# - Capitalize the genus names in column A (e.g. "justicia" -> "Justicia")
# - Rename the "image" header (E1) to "image url"
# - Populate an image-url value for every species row (E2:E6), replacing the
#   two old "stray" URLs that only covered a couple of rows
# - Drop the old hyperlink on E3 (it becomes a plain text/value cell, keeping
#   its existing "hyperlink" look style) since the whole column is now plain
#   URL text, not an embedded hyperlink object
# - Update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize genus names (column A)
$ws.Range("A2").Value = "Justicia"
$ws.Range("A3").Value = "Ruellia"
$ws.Range("A4").Value = "Acorus"
$ws.Range("A5").Value = "Sambucus"
$ws.Range("A6").Value = "Viburnum"

# Header rename: "image" -> "image url"
$ws.Range("E1").Value = "image url"

# Remove the hyperlink object on E3 (text stays, formatting stays)
$ws.Range("E3").Hyperlinks.Delete()

# Fill in image url column for every row
$ws.Range("E4").Value = "https://www.wrc.udel.edu/wp-content/heritage/viewtn.php?photo_id=56"
$ws.Range("E2").Value = "https://www.wrc.udel.edu/wp-content/heritage/viewtn.php?photo_id=1144"
$ws.Range("E5").Value = "https://www.wrc.udel.edu/wp-content/heritage/viewtn.php?photo_id=1556"
$ws.Range("E6").Value = "https://www.wrc.udel.edu/wp-content/heritage/viewtn.php?photo_id=1417"
$ws.Range("E3").Value = "https://www.wrc.udel.edu/wp-content/heritage/viewtn.php?photo_id=361"

# Move the active selection
$ws.Range("D9").Select() | Out-Null
